$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text edits (shared strings with multiple runs) ---
# A8: "Volume 31   Number  1" -> "...Number  2" (last character "1" -> "2")
$hdr1 = $ws.Range("A8")
$hdr1Len = $hdr1.Value2.Length
$hdr1.Characters($hdr1Len, 1).Text = "2"

# C9: "Report Covering the Week  1/1/2024  Through  1/7/2024"
#     -> "...1/8/2024  Through  1/14/2024"
$hdr2 = $ws.Range("C9")
$hdr2.Characters(27, 8).Text = "1/8/2024"
$hdr2.Characters(47, 8).Text = "1/14/2024"

# --- Crime statistics table edits (rows 15-27, 41, 43) ---
# Row 15
$ws.Range("D15").Copy($ws.Range("C15"))
$ws.Range("F15").Value = 3
$ws.Range("H27").Copy($ws.Range("N15"))
$ws.Range("N15").Value = 0
# Row 16
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 100
$ws.Range("I16").Value = 7
$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 133.333333333333
$ws.Range("L16").Value = -22.222222222222
$ws.Range("M16").Value = 133.333333333333
$ws.Range("N16").Value = -63.157894736842
# Row 17
$ws.Range("D15").Copy($ws.Range("C17"))
$ws.Range("E17").Value = -100
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -27.272727272727
$ws.Range("J17").Value = 4
$ws.Range("K17").Value = -75
$ws.Range("L17").Value = -80
$ws.Range("M17").Value = -83.333333333333
$ws.Range("N17").Value = -88.888888888888
# Row 18
$ws.Range("C18").Value = 4
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 33.333333333333
$ws.Range("F18").Value = 10
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = 66.666666666666
$ws.Range("I18").Value = 7
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 75
$ws.Range("L18").Value = 250
$ws.Range("M18").Value = 133.333333333333
$ws.Range("N18").Value = -22.222222222222
# Row 19
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 48
$ws.Range("G19").Value = 46
$ws.Range("H19").Value = 4.347826086956
$ws.Range("I19").Value = 22
$ws.Range("J19").Value = 21
$ws.Range("K19").Value = 4.761904761904
$ws.Range("L19").Value = -21.428571428571
$ws.Range("M19").Value = -18.518518518518
$ws.Range("N19").Value = -12
# Row 20
$ws.Range("C20").Value = 2
$ws.Range("D15").Copy($ws.Range("D20"))
$ws.Range("E15").Copy($ws.Range("E20"))
$ws.Range("F20").Value = 5
$ws.Range("H20").Value = 25
$ws.Range("I20").Value = 3
$ws.Range("K20").Value = 200
$ws.Range("L20").Value = -50
$ws.Range("N20").Value = -86.363636363636
# Row 21
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = 5.263157894736
$ws.Range("F21").Value = 84
$ws.Range("G21").Value = 72
$ws.Range("H21").Value = 16.666666666666
$ws.Range("I21").Value = 41
$ws.Range("J21").Value = 33
$ws.Range("K21").Value = 24.242424242424
$ws.Range("L21").Value = -18
$ws.Range("M21").Value = 5.128205128205
$ws.Range("N21").Value = -51.764705882352
# Row 22
$ws.Range("D15").Copy($ws.Range("D22"))
$ws.Range("E15").Copy($ws.Range("E22"))
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 0
$ws.Range("H27").Copy($ws.Range("L22"))
$ws.Range("L22").Value = -100
# Row 23
$ws.Range("D15").Copy($ws.Range("D23"))
$ws.Range("E15").Copy($ws.Range("E23"))
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 50
$ws.Range("L23").Value = -66.666666666666
# Row 24
$ws.Range("C24").Value = 15
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 50
$ws.Range("G24").Value = 37
$ws.Range("H24").Value = 35.135135135135
$ws.Range("I24").Value = 25
$ws.Range("J24").Value = 22
$ws.Range("K24").Value = 13.636363636363
$ws.Range("L24").Value = 31.578947368421
$ws.Range("M24").Value = 4.166666666666
# Row 25
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = 50
$ws.Range("F25").Value = 26
$ws.Range("G25").Value = 28
$ws.Range("H25").Value = -7.142857142857
$ws.Range("I25").Value = 16
$ws.Range("J25").Value = 14
$ws.Range("K25").Value = 14.285714285714
$ws.Range("L25").Value = 128.571428571429
$ws.Range("M25").Value = 6.666666666666
# Row 26
$ws.Range("D15").Copy($ws.Range("C26"))
$ws.Range("F26").Value = 3
# Row 27
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = -50
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("I27").Value = 2
$ws.Range("J27").Value = 3
$ws.Range("K27").Value = -33.333333333333
# Row 41
$ws.Range("J41").Value = 759
$ws.Range("K41").Value = 69.798657718120
$ws.Range("L41").Value = 31.770833333333
$ws.Range("M41").Value = -8.333333333333
$ws.Range("N41").Value = -34.342560553633
# Row 43
$ws.Range("J43").Value = 1276
$ws.Range("K43").Value = 37.5
$ws.Range("L43").Value = -23.638539796529
$ws.Range("M43").Value = -55.694444444444
$ws.Range("N43").Value = -69.260419176102
